$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.429.46'
$ws.Range("D3").Value = '3.173.12'
$ws.Range("E3").Value = '  -0.04%  '
$ws.Range("E4").Value = '  +0.10%  '
$__s = $ws.Range("D5").Style
$ws.Range("D5").Value = "'602.09"
$ws.Range("D5").Style = $__s
$ws.Range("E5").Value = '  +3.19%  '
$__s = $ws.Range("D6").Style
$ws.Range("D6").Value = "'135.95"
$ws.Range("D6").Style = $__s
$ws.Range("E6").Value = '  +0.85%  '
$ws.Range("E7").Value = '  -0.12%  '
$ws.Range("D8").Value = '3.172.59'
$ws.Range("E8").Value = '  +0.23%  '
$__s = $ws.Range("D9").Style
$ws.Range("D9").Value = "'0.514"
$ws.Range("D9").Style = $__s
$ws.Range("E9").Value = '  +2.54%  '
$ws.Range("E10").Value = '  +1.31%  '
$ws.Range("E11").Value = '  +2.68%  '
$__s = $ws.Range("D12").Style
$ws.Range("D12").Value = "'0.456"
$ws.Range("D12").Style = $__s
$ws.Range("E13").Value = '  +3.00%  '
$__s = $ws.Range("D14").Style
$ws.Range("D14").Value = "'34.87"
$ws.Range("D14").Style = $__s
$ws.Range("E14").Value = '  +5.42%  '
$ws.Range("D15").Value = '3.694.78'
$ws.Range("E15").Value = '  +0.15%  '
$ws.Range("E16").Value = '  +0.66%  '
$ws.Range("D17").Value = '3.172.86'
$ws.Range("E17").Value = '  -0.08%  '
$ws.Range("D18").Value = '63.419.94'
$ws.Range("E18").Value = '  +1.97%  '
$__s = $ws.Range("D19").Style
$ws.Range("D19").Value = "'6.60"
$ws.Range("D19").Style = $__s
$ws.Range("E19").Value = '  +0.83%  '
$__s = $ws.Range("D20").Style
$ws.Range("D20").Value = "'463.52"
$ws.Range("D20").Style = $__s
$ws.Range("E20").Value = '  +2.09%  '
$__s = $ws.Range("D21").Style
$ws.Range("D21").Value = "'14.00"
$ws.Range("D21").Style = $__s
$ws.Range("E21").Value = '  +1.02%  '
$ws.Range("E22").Value = '  -0.38%  '
$__s = $ws.Range("D23").Style
$ws.Range("D23").Value = "'7.69"
$ws.Range("D23").Style = $__s
$ws.Range("E23").Value = '  +1.39%  '
$__s = $ws.Range("D24").Style
$ws.Range("D24").Value = "'13.31"
$ws.Range("D24").Style = $__s
$ws.Range("E24").Value = '  +0.45%  '
$__s = $ws.Range("D25").Style
$ws.Range("D25").Value = "'83.24"
$ws.Range("D25").Style = $__s
$ws.Range("E25").Value = '  +1.39%  '
$ws.Range("E26").Value = '  +0.13%  '
$__s = $ws.Range("D27").Style
$ws.Range("D27").Value = "'2.71"
$ws.Range("D27").Style = $__s
$ws.Range("E27").Value = '  +1.22%  '
$ws.Range("E28").Value = '  +0.05%  '
$__s = $ws.Range("D29").Style
$ws.Range("D29").Value = "'2.09"
$ws.Range("D29").Style = $__s
$ws.Range("E29").Value = '  +4.30%  '
$ws.Range("E30").Value = '  -1.01%  '
$ws.Range("E31").Value = '  -0.67%  '
$__s = $ws.Range("D32").Style
$ws.Range("D32").Value = "'27.19"
$ws.Range("D32").Style = $__s
$ws.Range("E32").Value = '  +0.29%  '
$ws.Range("E33").Value = '  -0.83%  '
$__s = $ws.Range("D34").Style
$ws.Range("D34").Value = "'2.43"
$ws.Range("D34").Style = $__s
$ws.Range("E34").Value = '  +2.14%  '
$ws.Range("E35").Value = '  -1.48%  '
$__s = $ws.Range("D36").Style
$ws.Range("D36").Value = "'5.91"
$ws.Range("D36").Style = $__s
$ws.Range("E36").Value = '  +2.40%  '
$ws.Range("D37").Value = '0.0₃0734'
$ws.Range("E37").Value = '  +6.94%  '
$__s = $ws.Range("D38").Style
$ws.Range("D38").Value = "'51.25"
$ws.Range("D38").Style = $__s
$ws.Range("E38").Value = '  +0.47%  '
$__s = $ws.Range("D39").Style
$ws.Range("D39").Value = "'0.0391"
$ws.Range("D39").Style = $__s
$ws.Range("E39").Value = '  +1.64%  '
$ws.Range("E40").Value = '  +2.01%  '
$ws.Range("E41").Value = '  +1.17%  '
$ws.Range("E42").Value = '  +0.33%  '
$__s = $ws.Range("D43").Style
$ws.Range("D43").Value = "'394.91"
$ws.Range("D43").Style = $__s
$ws.Range("E43").Value = '  -3.15%  '
$ws.Range("D44").Value = '2.807.63'
$ws.Range("E44").Value = '  -4.63%  '
$ws.Range("E45").Value = '  +1.15%  '
$__s = $ws.Range("D46").Style
$ws.Range("D46").Value = "'36.24"
$ws.Range("D46").Style = $__s
$ws.Range("E46").Value = '  +1.63%  '
$ws.Range("B47").Value = 'Fetch.AI'
$ws.Range("C47").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$__s = $ws.Range("D47").Style
$ws.Range("D47").Value = "'2.12"
$ws.Range("D47").Style = $__s
$ws.Range("E47").Value = '  -0.29%  '
$ws.Range("B48").Value = 'USDe'
$ws.Range("C48").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$__s = $ws.Range("D48").Style
$ws.Range("D48").Value = "'0.999"
$ws.Range("D48").Style = $__s
$ws.Range("E48").Value = '  +0.01%  '
$__s = $ws.Range("D49").Style
$ws.Range("D49").Value = "'126.12"
$ws.Range("D49").Style = $__s
$ws.Range("E49").Value = '  +2.34%  '
$__s = $ws.Range("D50").Style
$ws.Range("D50").Value = "'25.21"
$ws.Range("D50").Style = $__s
$ws.Range("E50").Value = '  -0.47%  '
$ws.Range("E51").Value = '  +1.00%  '
